$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 11; $row++) {
    $ws.Cells.Item($row, 1).Value = "sku"
    $ws.Cells.Item($row, 2).Value = "name"
    $ws.Cells.Item($row, 3).Value = "quantity"
    $ws.Cells.Item($row, 4).Value = "cost_per"
    $ws.Cells.Item($row, 5).Value = "total_cost"
}
